$d = $word.ActiveDocument

# --- Step 0: grab the "Meta description" bold run (paragraph 2) before anything
#     shifts paragraph indices. This run is <w:rPr><w:b/></w:rPr><w:t>Meta description</w:t>.
$metaPara = $d.Paragraphs.Item(2)
$metaStart = $metaPara.Range.Start
$boldSrc = $d.Range($metaStart, $metaStart + 16)

# --- Step 1: insert a brand-new paragraph right before the final
#     ("Please create an image...") paragraph, carrying the bold run text/formatting,
#     then retarget the text to the new heading copy.
$count = $d.Paragraphs.Count
$secondLast = $d.Paragraphs.Item($count - 1)
$secondLast.Range.InsertParagraphAfter()

$newP = $d.Paragraphs.Item($count)
$newP.Style = $d.Styles.Item("Normal")
$newRange = $newP.Range
$newRange.FormattedText = $boldSrc.FormattedText

$newHeadStart = $newP.Range.Start
$boldDest = $d.Range($newHeadStart, $newHeadStart + 17)
$boldDest.Text = "Play Diamond Duke for Free - Classic Bar Slot Machine Game"

# --- Step 2: replace the text of the (now shifted) final paragraph - the italic
#     image-prompt paragraph - with the old meta-description body text.
$count = $d.Paragraphs.Count
$lastP = $d.Paragraphs.Item($count)
$lastRange = $lastP.Range
$lastTextRange = $d.Range($lastRange.Start, $lastRange.End - 1)
$lastTextRange.Text = "Read our review of Diamond Duke, a classic bar slot machine game with interesting win multipliers and special symbols. Play it for free today!"

# --- Step 3: delete the original "Meta description" paragraph entirely.
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

Write-Output "done"
